$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 9
$ws.Range("G9").Value = 2.05
$ws.Range("H9").Value = 3.5
$ws.Range("I9").Value = 3.4
$ws.Range("J9").Value = 2.75
$ws.Range("L9").Value = 4
$ws.Range("N9").Value = 9.5
$ws.Range("S9").Value = 1.44
$ws.Range("T9").Value = 2.63
$ws.Range("X9").Value = 9.5
$ws.Range("Z9").Value = 19
$ws.Range("AD9").Value = 6.5
$ws.Range("AE9").Value = 15
$ws.Range("AH9").Value = 9.5
$ws.Range("AI9").Value = 17
$ws.Range("AJ9").Value = 12
$ws.Range("AL9").Value = 29
$ws.Range("AO9").Value = 11
$ws.Range("AQ9").Value = 41
$ws.Range("AT9").Value = 2.63
$ws.Range("AU9").Value = 8
$ws.Range("AX9").Value = 5.5
$ws.Range("BA9").Value = 67
$ws.Range("BB9").Value = 81
$ws.Range("BC9").Value = 201

# Row 11
$ws.Range("G11").Value = 2.25
$ws.Range("H11").Value = 2.9
$ws.Range("I11").Value = 3.6
$ws.Range("J11").Value = 3.2
$ws.Range("L11").Value = 4.5
$ws.Range("M11").Value = 1.14
$ws.Range("N11").Value = 5.5
$ws.Range("O11").Value = 1.67
$ws.Range("P11").Value = 2.1
$ws.Range("Q11").Value = 3.4
$ws.Range("R11").Value = 1.33
$ws.Range("X11").Value = 9
$ws.Range("AA11").Value = 26
$ws.Range("AC11").Value = 5
$ws.Range("AI11").Value = 15
$ws.Range("AO11").Value = 15

# Row 12
$ws.Range("G12").Value = 3.5
$ws.Range("I12").Value = 2.2
$ws.Range("J12").Value = 4
$ws.Range("K12").Value = 2
$ws.Range("L12").Value = 3
$ws.Range("U12").Value = 2
$ws.Range("V12").Value = 1.73
$ws.Range("W12").Value = 8.5
$ws.Range("Y12").Value = 13
$ws.Range("AH12").Value = 6.5
$ws.Range("AI12").Value = 9.5
$ws.Range("AK12").Value = 21
$ws.Range("AN12").Value = 5
$ws.Range("AQ12").Value = 67
$ws.Range("AR12").Value = 101
$ws.Range("AY12").Value = 13
